# Update fsl_flag_description file (FSL sheet) - refine rationale/action text
# to add explicit thresholds, matching the author's clarified wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FSL")

$ws.Range("C7").Value = 'If reported cereal consumption is low (cereal<5 days), review according to food sources and livelihoods sytems for coherence. '
$ws.Range("C8").Value = 'If reported oil consumption is low (oill<5 days), review according to food sources and livelihoods sytems for coherence. '
$ws.Range("B9").Value = 'Except if in EXTREME cases (e.g. Somalia), it will be very rare for many/any HHs to have such low FCS scores (FCS <10)'
$ws.Range("B10").Value = 'FCS total score is 112 which corresponds to consuming all food groups every day, highly unlikely in any context. WFP considered FCS>56 already as high. To be reviewed against context for coherence.'
$ws.Range("B12").Value = ' If fcs_score < 35 and rcsi_score <= 4. FCS and rCSI are strongly correlated with a negative relation. This combination is''t impossible but needs to be verified within its context. If collected in early stages of a crisis, would be expected that rCSI will still be sensitive enough to measure early consumption-based coping and measure a higher score.'
$ws.Range("B15").Value = "If rCSI score is high (in contexts of early crisis stages) while protein consumption is also reported as frequent (min. 5 days), it could be that either dimension wasn't understood or collected properly. "
$ws.Range("B20").Value = 'If HHS score >= 5. Reports of severe and very severe hunger are flags for particularly At Risk Households suffering from prolonged acute food insecurity '
$ws.Range("B32").Value = 'If HH has both rCSI score>18 and FCS score>56. Any HH that would have an acceptable FCS score (higher scores) and a high rCSI score is most likely indicative of data quality issue with one or both indicators'

# Apply a plain "Normal" style across the used columns (matches the author's
# whole-column formatting pass visible in the workbook XML).
$ws.Columns("A:C").Font.Name = "Calibri"
$ws.Columns("A:C").Font.Size = 11

# Restore the view to the top of the sheet with B7 as the active cell.
$ws.Range("B7").Select
